# Update the dSF column (F) values for rows 2-18, 20-21 as part of a
# repull/recalculation of the mean delta-S (final) statistic.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 4
    3  = 4
    4  = 3
    5  = 3
    6  = -4
    7  = -1
    8  = -5
    9  = -2
    10 = 1
    11 = -2
    12 = 3
    13 = -2
    14 = 0
    15 = -1
    16 = 2
    17 = 1
    18 = -4
    20 = -1
    21 = 1
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
